$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 0.9943347821518321
$ws.Range("C3").Value = 0.99393435282297
$ws.Range("D3").Value = 0.9940335043889796

$ws.Range("B5").Value = 0.9850277839156065
$ws.Range("C5").Value = 0.984151486532484
$ws.Range("D5").Value = 0.9843573581379104
